$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 11
    4  = 4
    5  = 7
    6  = 8
    7  = 3
    8  = 6
    9  = 3
    10 = 10
    11 = 7
    12 = 3
    13 = 10
    14 = 6
    15 = 8
    16 = 3
    17 = 3
    18 = 6
    19 = 7
    20 = 6
    21 = 1
    22 = 7
    23 = 8
    24 = 8
    25 = 7
    26 = 7
    27 = 5
    28 = 4
    29 = 3
    30 = 9
    31 = 6
    32 = 9
    33 = 5
    34 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
